$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 2024 participant refresh ---
# The old dataset (14 participant rows, rows 2-15) is replaced in full by an
# updated export: some institutions removed (DSS, Univ. i Bergen), several added
# (Innlandet, Molde, Ostfold, Vestlandet, Kristiania, NTNU/St. Olavs, UiT/UNN,
# Sorost-Norge), Gjovik + Sor-Trondelag folded into the NTNU row, and the
# agreement window + label bumped from "Emerald 2015" (2015-01-01..2015-12-31) to
# "Emerald 2023-2025" (2023-01-01..2025-12-31), with agr_date_updated -> 2024-03-18.
# The new table needs 19 data rows (2-20), five more than the old 14 (2-15), so
# first stretch the date formatting (columns I:K) from the last existing row down
# over the new rows before the values are (re)written.
$ws.Range("I15:K15").Copy($ws.Range("I16:K20"))

# Clear out the previous data rows before laying down the refreshed table.
$ws.Range("A2:K20").ClearContents()

$data = @(
  @('emerald', 'Handelshøyskolen BI', 'bi', 'Handelshøyskolen BI', 158, 15353, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'Høgskolen I Innlandet', 'inn', 'Høgskolen i Innlandet', 209, 15428, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'Høgskolen i Molde - Vitenskapelig høgskole i logistikk', 'himolde', 'Høgskolen i Molde - Vitenskapelig høgskole i logistikk', 211, 15435, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'Høgskolen i Østfold', 'hioef', 'Høgskolen i Østfold', 224, 15443, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'Høgskulen På Vestlandet', 'hvl', 'Høgskulen på Vestlandet', 203, 15449, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'Høyskolen Kristiania - Ernst G Mortensens Stiftelse', 'hk', 'Høyskolen Kristiania', 1615, 15444, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'Nasjonalbiblioteket', 'nb', 'Nasjonalbiblioteket', 5931, 15471, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'Nofima', 'nofima', 'NOFIMA', 7543, 15496, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'Nord Universitet', 'nord', 'Nord universitet', 204, 15499, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'Norges Handelshøyskole', 'nhh', 'Norges Handelshøyskole', 191, 15480, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'Norges teknisk-naturvitenskapelige universitet', 'ntnu', 'Norges teknisk-naturvitenskapelige universitet', 194, 15507, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'Norges teknisk-naturvitenskapelige universitet', 'ntnu_1920', 'St. Olavs Hospital HF', 1920, 15507, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'OsloMet - Storbyuniversitetet', 'oslomet', 'OsloMet - storbyuniversitetet', 215, 15439, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'UiO : Universitetsbiblioteket', 'uio', 'Universitetet i Oslo', 185, 15550, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'UiT - Norges arktiske universitet', 'uit', 'UiT Norges arktiske universitet', 186, 15552, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'UiT - Norges arktiske universitet', 'uit_1902', 'Universitetssykehuset Nord-Norge HF', 1902, 15552, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'Universitetet I Stavanger', 'uis', 'Universitetet i Stavanger', 217, 15551, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'Universitetet I Sørøst-Norge', 'usn', 'Universitetet i Sørøst-Norge', 222, 15447, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369),
  @('emerald', 'Universitetet i Agder', 'uia', 'Universitetet i Agder', 201, 15548, 'Emerald Publishing', 'Emerald 2023-2025', 44927, 46022, 45369)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
